$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-20 Wednesday" "2025-08-21 Thursday"

Replace-Text "447×3=" "374×9="
Replace-Text "920×8=" "190×7="
Replace-Text "917×2=" "681×9="
Replace-Text "270×8=" "260×5="
Replace-Text "788×5=" "775×6="

Replace-Text "559×5=" "224×4="
Replace-Text "323×8=" "103×7="
Replace-Text "741×8=" "610×2="
Replace-Text "230×6=" "459×5="
Replace-Text "249×4=" "697×8="

Replace-Text "985×5=" "437×7="
Replace-Text "737×4=" "445×9="
Replace-Text "381×5=" "264×7="
Replace-Text "889×6=" "726×9="
Replace-Text "272×6=" "550×5="

Replace-Text "790×7=" "659×2="
Replace-Text "319×4=" "844×6="
Replace-Text "595×7=" "280×6="
Replace-Text "494×9=" "766×3="
Replace-Text "529×8=" "579×5="

Replace-Text "799×7=" "689×3="
Replace-Text "504×4=" "326×3="
Replace-Text "455×5=" "545×8="
Replace-Text "315×3=" "200×3="
Replace-Text "397×2=" "950×8="
